# Scene 4A edits — Mom's morning-kitchen dialogue pass ("write some new for stephen")
#
# Strategy: locate each target paragraph by its (unique) current text via
# Paragraphs.Item(n).Range, then set .Text directly (which is what gives the
# expected xml:space="preserve" behaviour on the rewritten runs). New lines of
# dialogue are spliced in with Range.InsertParagraphAfter() immediately before
# setting their text, so they pick up the same paragraph/run formatting as
# their neighbours.
#
# We walk the document from the bottom up so earlier insertions never shift
# the paragraph indices we still need to touch.

$d = $word.ActiveDocument

function Get-ParaIndexByText($needle) {
    $paras = $d.Paragraphs
    $n = $paras.Count
    for ($i = $n; $i -ge 1; $i--) {
        # Paragraph.Range.Text includes the trailing paragraph mark (\r), so
        # trim before comparing against a plain-text needle.
        if ($paras.Item($i).Range.Text.TrimEnd() -eq $needle) {
            return $i
        }
    }
    throw "Get-ParaIndexByText: no paragraph matched [$needle]"
}

# --- 12) "Mom (neutral smiling): I'm going to work now..." -> "Mom (waving smiling): ..."
$idx = Get-ParaIndexByText("Mom (neutral smiling): I’m going to work now, so I’ll see you later, okay?")
$d.Paragraphs.Item($idx).Range.Text = "Mom (waving smiling): I’m going to work now, so I’ll see you later, okay?"

# --- 11) "Mom (neutral smiling): Thanks again, Pro." -> "Mom: Thanks again, Pro."
$idx = Get-ParaIndexByText("Mom (neutral smiling): Thanks again, Pro.")
$d.Paragraphs.Item($idx).Range.Text = "Mom: Thanks again, Pro."

# --- 10) "Mom (neutral smiling_nervous): I know. Do your best though, okay?" -> "Mom (neutral smiling): ..."
$idx = Get-ParaIndexByText("Mom (neutral smiling_nervous): I know. Do your best though, okay?")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral smiling): I know. Do your best though, okay?"

# --- 9) "Mom (neutral smiling_nervous): Everything going alright in school?" -> "Mom (neutral curious): ..."
$idx = Get-ParaIndexByText("Mom (neutral smiling_nervous): Everything going alright in school?")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral curious): Everything going alright in school?"

# --- 8) Merge the 3-run "How have you been? / You / been doing okay recently?" paragraph
#        into a single run, and swap "neutral smiling" for "neutral worried_slightly".
$idx = Get-ParaIndexByText("Mom (neutral smiling): How have you been? You been doing okay recently?")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral worried_slightly): How have you been? You been doing okay recently?"

# --- 7) "Mom (neutral smiling): It's good to see that you're starting to eat a little more." -> smiling_eyes_closed
$idx = Get-ParaIndexByText("Mom (neutral smiling): It’s good to see that you’re starting to eat a little more.")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral smiling_eyes_closed): It’s good to see that you’re starting to eat a little more."

# --- 6) "Pro: No, I'm not. Eggs are always good." -> "Pro: Nope. Eggs are always good."
$idx = Get-ParaIndexByText("Pro: No, I’m not. Eggs are always good.")
$d.Paragraphs.Item($idx).Range.Text = "Pro: Nope. Eggs are always good."

# --- 5) "Mom (neutral worried): Are you tired of them?" -> "Mom (neutral curious): ..."
$idx = Get-ParaIndexByText("Mom (neutral worried): Are you tired of them?")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral curious): Are you tired of them?"

# --- 4) "Mom (neutral smiling): A couple weeks ago..." -> smiling_eyes_closed
$idx = Get-ParaIndexByText("Mom (neutral smiling): A couple weeks ago they were on sale, so I bought a lot. But they’re gonna expire soon, so I wanted us to eat them all so we don’t have to throw any away.")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral smiling_eyes_closed): A couple weeks ago they were on sale, so I bought a lot. But they’re gonna expire soon, so I wanted us to eat them all so we don’t have to throw any away."

# --- 3) New paragraph "Mom (neutral smiling):" inserted right before
#        "A couple minutes later she finishes up..."
$idx = Get-ParaIndexByText("I sit down as instructed and put my head down on the table, watching as my mom puts the finishing touches on our meal. It’s amazing how easy she makes cooking look, especially at this time. If I tried to cook this early in the morning, I’d probably fall asleep and set the house on fire.")
$anchor = $d.Paragraphs.Item($idx).Range
$anchor.InsertParagraphAfter()
$d.Paragraphs.Item($idx + 1).Range.Text = "Mom (neutral smiling):"

# --- 2) New paragraph "Mom (exit):" inserted right after
#        "Mom (neutral smiling): Thanks. Have a seat, I'll be done with breakfast in a few minutes."
$idx = Get-ParaIndexByText("Mom (neutral smiling): Thanks. Have a seat, I’ll be done with breakfast in a few minutes.")
$anchor = $d.Paragraphs.Item($idx).Range
$anchor.InsertParagraphAfter()
$d.Paragraphs.Item($idx + 1).Range.Text = "Mom (exit):"
# Re-stamp the (unchanged-text) "Thanks. Have a seat..." run so it picks up the
# xml:space="preserve" the target diff marks on it.
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral smiling): Thanks. Have a seat, I’ll be done with breakfast in a few minutes."

# --- 1) "Mom (neutral smiling): I have to go to work early today..." -> "Mom (neutral neutral): ..."
$idx = Get-ParaIndexByText("Mom (neutral smiling): I have to go to work early today, so since you’re up could you wash the dishes when you’re done eating?")
$d.Paragraphs.Item($idx).Range.Text = "Mom (neutral neutral): I have to go to work early today, so since you’re up could you wash the dishes when you’re done eating?"

Write-Output "done"
